$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.161.97'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.989.17'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '501.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.37%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  -3.31%  '
$ws.Range('E9').Value = '  -4.07%  '
$ws.Range('E10').Value = '  -4.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.359'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('D12').Value = '3.504.82'
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.17'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.75%  '
$ws.Range('E15').Value = '  -5.64%  '
$ws.Range('D16').Value = '57.186.57'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '2.987.71'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.66'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.94%  '
$ws.Range('E20').Value = '  -3.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.30%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.491'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.32%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.163'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('D28').Value = '0.0₃0897'
$ws.Range('E28').Value = '  -7.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.56'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('E31').Value = '  -3.85%  '
$ws.Range('E32').Value = '  -4.80%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '155.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('E37').Value = '  -5.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.97%  '
$ws.Range('E39').Value = '  -5.84%  '
$ws.Range('D40').Value = '3.020.90'
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('E45').Value = '  -5.94%  '
$ws.Range('D46').Value = '2.197.21'
$ws.Range('E46').Value = '  -5.71%  '
$ws.Range('B47').Value = 'Cosmos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.95'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.52%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.941'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -8.79%  '
$ws.Range('E49').Value = '  -4.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.37%  '
$ws.Range('E51').Value = '  -11.08%  '
